# Updated cryptos list (price/volume refresh) matching the automated
# "Updated cryptos list ... with GitHub Actions" commit.
#
# Price (column D) and Volume(1h) (column E) are stored as plain text in
# the sheet (e.g. "44.099.78", "  +1.98%  "), so a leading apostrophe is
# used for values that would otherwise be auto-coerced to a number by
# Excel's Value setter (single-dot decimals like "238.62").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "44.099.78"
$ws.Range("E2").Value = "  +1.98%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.358.05"
$ws.Range("E3").Value = "  +0.48%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.10%  "

# Row 5 - XRP
$ws.Range("E5").Value = "  +3.98%  "

# Row 6 - BNB
$ws.Range("D6").Value = "'238.62"
$ws.Range("E6").Value = "  +3.32%  "

# Row 7 - Solana
$ws.Range("D7").Value = "'72.95"
$ws.Range("E7").Value = "  +11.78%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "'0.543"
$ws.Range("E9").Value = "  +18.68%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +6.29%  "

# Row 11 - Avalanche
$ws.Range("D11").Value = "'29.35"
$ws.Range("E11").Value = "  +10.30%  "

# Row 12 - TRON
$ws.Range("D12").Value = "'0.108"
$ws.Range("E12").Value = "  +2.93%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "2.713.23"
$ws.Range("E13").Value = "  +0.64%  "

# Row 14 - Chainlink
$ws.Range("E14").Value = "  +10.02%  "

# Row 15 - Polkadot
$ws.Range("D15").Value = "'6.71"
$ws.Range("E15").Value = "  +7.56%  "

# Row 16 - Polygon
$ws.Range("D16").Value = "'0.907"
$ws.Range("E16").Value = "  +8.25%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.361.64"
$ws.Range("E17").Value = "  +0.70%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "44.060.02"
$ws.Range("E18").Value = "  +1.83%  "

# Row 19 - ShibaInu
$ws.Range("E19").Value = "  +5.13%  "

# Row 20 - Litecoin
$ws.Range("E20").Value = "  +6.09%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  +4.67%  "

# Row 22 - BitcoinCash
$ws.Range("D22").Value = "'255.62"
$ws.Range("E22").Value = "  +3.10%  "

# Row 24 - WEMIXToken
$ws.Range("E24").Value = "  -3.79%  "

# Row 25 - PancakeSwap
$ws.Range("E25").Value = "  +3.76%  "

# Row 26 - Cosmos
$ws.Range("D26").Value = "'10.50"
$ws.Range("E26").Value = "  +6.48%  "

# Row 27 - Toncoin
$ws.Range("E27").Value = "  -1.36%  "

# Row 28 - EthereumClassic
$ws.Range("D28").Value = "'22.44"
$ws.Range("E28").Value = "  +1.15%  "

# Row 29/30 swap: ImmutableX <-> Monero traded ranking places
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "'173.04"
$ws.Range("E29").Value = "  -1.27%  "

$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").Value = "'1.58"
$ws.Range("E30").Value = "  +6.17%  "

# Row 31 - Kaspa
$ws.Range("D31").Value = "'0.133"
$ws.Range("E31").Value = "  +2.52%  "

# Row 32 - Stellar
$ws.Range("E32").Value = "  +6.06%  "

# Row 33 - Filecoin
$ws.Range("D33").Value = "'5.21"
$ws.Range("E33").Value = "  +4.77%  "

# Row 34 - Hedera
$ws.Range("D34").Value = "'0.0733"
$ws.Range("E34").Value = "  +6.97%  "

# Row 35 - InternetComputer(DFINITY)
$ws.Range("D35").Value = "'5.26"
$ws.Range("E35").Value = "  +5.59%  "

# Row 36 - RenderToken
$ws.Range("D36").Value = "'3.93"
$ws.Range("E36").Value = "  +10.07%  "

# Row 37 - LidoDAOToken
$ws.Range("E37").Value = "  -1.06%  "

# Row 38 - THORChain
$ws.Range("E38").Value = "  -0.02%  "

# Row 39 - VeChain
$ws.Range("D39").Value = "'0.0271"
$ws.Range("E39").Value = "  +7.25%  "

# Row 40 - InjectiveProtocol
$ws.Range("D40").Value = "'19.60"
$ws.Range("E40").Value = "  +9.84%  "

# Row 41 - BinanceUSD
$ws.Range("E41").Value = "  +0.06%  "

# Row 42 - FraxShare
$ws.Range("D42").Value = "'8.88"
$ws.Range("E42").Value = "  +0.00%  "

# Row 43 - TrustWalletToken
$ws.Range("E43").Value = "  +4.12%  "

# Row 44 - Cronos
$ws.Range("E44").Value = "  +4.17%  "

# Row 45 - ARBITRUM
$ws.Range("E45").Value = "  +1.43%  "

# Row 46 - Aave
$ws.Range("D46").Value = "'98.66"
$ws.Range("E46").Value = "  +0.32%  "

# Row 47 - FTXToken
$ws.Range("E47").Value = "  +1.95%  "

# Row 48 - Algorand
$ws.Range("E48").Value = "  +12.59%  "

# Row 49 - NEARProtocol
$ws.Range("E49").Value = "  +5.26%  "

# Row 50 - Maker
$ws.Range("D50").Value = "1.441.61"
$ws.Range("E50").Value = "  +0.51%  "

# Row 51 - HuobiToken
$ws.Range("E51").Value = "  +1.50%  "
